# Update NATMI TPM-derived values for the Il4-Cd53 LR pair sheet.
# The underlying analysis was re-run with new TPM values, which changes
# several per-row derived statistics (detection counts/rates, expression
# values, specificities, and edge weights) for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Il4/Cd53 -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3454506666666666
$ws.Range("H2").Value = 1.036352
$ws.Range("I2").Value = 0.1052716477644991
$ws.Range("J2").Value = 0.1052716477644991
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.051095
$ws.Range("N2").Value = 0.153285
$ws.Range("Q2").Value = 0.01765080181333333
$ws.Range("R2").Value = 0.15885721632
$ws.Range("S2").Value = 0.1052716477644991
$ws.Range("T2").Value = 0.1052716477644991

# Row 3 (FAPs -> Il4/Cd53 -> ECs)
$ws.Range("I3").Value = 0.3398937483175971
$ws.Range("J3").Value = 0.3398937483175971
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.051095
$ws.Range("N3").Value = 0.153285
$ws.Range("Q3").Value = 0.056989676865
$ws.Range("R3").Value = 0.512907091785
$ws.Range("S3").Value = 0.3398937483175971
$ws.Range("T3").Value = 0.3398937483175971

# Row 4 (MuSCs -> Il4/Cd53 -> ECs)
$ws.Range("I4").Value = 0.5548346039179038
$ws.Range("J4").Value = 0.5548346039179038
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.051095
$ws.Range("N4").Value = 0.153285
$ws.Range("Q4").Value = 0.093028615405
$ws.Range("R4").Value = 0.8372575386450001
$ws.Range("S4").Value = 0.5548346039179038
$ws.Range("T4").Value = 0.5548346039179038
